# Auto-generated Excel COM-interop edit script
# Applies cell-value corrections to the Leviathan Profits (Leve) tracker sheets
# as produced by the scheduled data-refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -54992

$ws.Range("H96").Value = 334708.5
$ws.Range("I96").Value = 400804.4
$ws.Range("J96").Value = 4229
$ws.Range("K96").Value = 1202413.2
$ws.Range("L96").Value = 12687
$ws.Range("M96").Value = -1201040.2
$ws.Range("N96").Value = -15433

$ws.Range("H100").Value = 3262
$ws.Range("I100").Value = 2699.4443
$ws.Range("J100").Value = 4274.6
$ws.Range("K100").Value = 2699.4443
$ws.Range("L100").Value = 4274.6
$ws.Range("M100").Value = -2158.4443
$ws.Range("N100").Value = -5356.6

$ws.Range("H106").Value = 124886
$ws.Range("I106").Value = 169164.83
$ws.Range("J106").Value = 36328.332
$ws.Range("K106").Value = 169164.83
$ws.Range("L106").Value = 36328.332
$ws.Range("M106").Value = -168533.83
$ws.Range("N106").Value = -37590.332

$ws.Range("H127").Value = 74705.84
$ws.Range("I127").Value = 80431.336
$ws.Range("K127").Value = 241294.008
$ws.Range("M127").Value = -236334.008

$ws.Range("H132").Value = 2809.2144
$ws.Range("I132").Value = 982.67645
$ws.Range("K132").Value = 2948.02935
$ws.Range("M132").Value = -418.0293500000002

$ws.Range("H135").Value = 56207.723
$ws.Range("I135").Value = 613.0909
$ws.Range("K135").Value = 5517.8181
$ws.Range("M135").Value = -2982.8181

$ws.Range("H138").Value = 2751
$ws.Range("I138").Value = 1351.4
$ws.Range("K138").Value = 4054.2
$ws.Range("M138").Value = 1085.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 25669.5
$ws.Range("I16").Value = 254.5
$ws.Range("J16").Value = 76499.5
$ws.Range("K16").Value = 254.5
$ws.Range("L16").Value = 76499.5
$ws.Range("M16").Value = 32.5
$ws.Range("N16").Value = -77073.5

$ws.Range("H32").Value = 156973.17
$ws.Range("I32").Value = 172754.34
$ws.Range("K32").Value = 172754.34
$ws.Range("M32").Value = -172467.34

$ws.Range("H102").Value = 3666.5557
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 10631.782
$ws.Range("I122").Value = 11860.15
$ws.Range("K122").Value = 35580.45
$ws.Range("M122").Value = -33130.45

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1748.2858
$ws.Range("I86").Value = 1532
$ws.Range("J86").Value = 2137.6
$ws.Range("K86").Value = 1532
$ws.Range("L86").Value = 2137.6
$ws.Range("M86").Value = -409
$ws.Range("N86").Value = -4383.6

$ws.Range("H89").Value = 1748.2858
$ws.Range("I89").Value = 1532
$ws.Range("J89").Value = 2137.6
$ws.Range("K89").Value = 7660
$ws.Range("L89").Value = 10688
$ws.Range("M89").Value = -2044
$ws.Range("N89").Value = -21920

$ws.Range("H94").Value = 1249.6364
$ws.Range("I94").Value = 1056.8
$ws.Range("J94").Value = 1410.3334
$ws.Range("K94").Value = 1056.8
$ws.Range("L94").Value = 1410.3334
$ws.Range("M94").Value = -605.8
$ws.Range("N94").Value = -2312.3334

$ws.Range("H105").Value = 5004166
$ws.Range("J105").Value = 1704.2
$ws.Range("L105").Value = 1704.2
$ws.Range("N105").Value = -5198.2

$ws.Range("H107").Value = 10164.286
$ws.Range("I107").Value = 3450.2942
$ws.Range("K107").Value = 3450.2942
$ws.Range("M107").Value = -1530.2942

$ws.Range("H134").Value = 2236
$ws.Range("I134").Value = 2236
$ws.Range("K134").Value = 6708
$ws.Range("M134").Value = -4173

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3385.6667
$ws.Range("I31").Value = 3162
$ws.Range("J31").Value = 3833
$ws.Range("K31").Value = 3162
$ws.Range("L31").Value = 3833
$ws.Range("M31").Value = -2867
$ws.Range("N31").Value = -4423

$ws.Range("H34").Value = 3385.6667
$ws.Range("I34").Value = 3162
$ws.Range("J34").Value = 3833
$ws.Range("K34").Value = 3162
$ws.Range("L34").Value = 3833
$ws.Range("M34").Value = -2960
$ws.Range("N34").Value = -4237

$ws.Range("H99").Value = 3297
$ws.Range("I99").Value = 3163.6667
$ws.Range("J99").Value = 3497
$ws.Range("K99").Value = 3163.6667
$ws.Range("L99").Value = 3497
$ws.Range("M99").Value = -1665.6667
$ws.Range("N99").Value = -6493

$ws.Range("H126").Value = 3297
$ws.Range("I126").Value = 3163.6667
$ws.Range("J126").Value = 3497
$ws.Range("K126").Value = 9491.000100000001
$ws.Range("L126").Value = 10491
$ws.Range("M126").Value = -7021.000100000001
$ws.Range("N126").Value = -15431

$ws.Range("H132").Value = 1722.9166
$ws.Range("I132").Value = 1722.9166
$ws.Range("K132").Value = 5168.7498
$ws.Range("M132").Value = -2638.7498

$ws.Range("H134").Value = 3131.8635
$ws.Range("I134").Value = 3024
$ws.Range("K134").Value = 9072
$ws.Range("M134").Value = -6537

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 4779.8
$ws.Range("I93").Value = 800
$ws.Range("K93").Value = 2400
$ws.Range("M93").Value = -528

$ws.Range("H118").Value = 4077.2
$ws.Range("I118").Value = 200
$ws.Range("J118").Value = 4354.143
$ws.Range("K118").Value = 600
$ws.Range("L118").Value = 13062.429
$ws.Range("M118").Value = 643
$ws.Range("N118").Value = -15548.429

$ws.Range("H128").Value = 325911.44
$ws.Range("I128").Value = 325911.44
$ws.Range("K128").Value = 977734.3200000001
$ws.Range("M128").Value = -972754.3200000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11776.929
$ws.Range("I80").Value = 22497.8
$ws.Range("K80").Value = 22497.8
$ws.Range("M80").Value = -21499.8

$ws.Range("H83").Value = 11776.929
$ws.Range("I83").Value = 22497.8
$ws.Range("K83").Value = 112489
$ws.Range("M83").Value = -107497

$ws.Range("H102").Value = 3409.5833
$ws.Range("I102").Value = 4130.7144
$ws.Range("J102").Value = 2400
$ws.Range("K102").Value = 4130.7144
$ws.Range("L102").Value = 2400
$ws.Range("M102").Value = -2508.7144
$ws.Range("N102").Value = -5644

$ws.Range("H105").Value = 96760
$ws.Range("J105").Value = 96760
$ws.Range("L105").Value = 96760
$ws.Range("N105").Value = -103748

$ws.Range("H132").Value = 870.5714
$ws.Range("I132").Value = 794.4
$ws.Range("J132").Value = 1061
$ws.Range("K132").Value = 2383.2
$ws.Range("L132").Value = 3183
$ws.Range("M132").Value = 146.8000000000002
$ws.Range("N132").Value = -8243

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 68499.75
$ws.Range("I7").Value = 68499.75
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 68499.75
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -68387.75
$ws.Range("N7").ClearContents()

$ws.Range("H22").Value = 1418.0588
$ws.Range("I22").Value = 2295.2
$ws.Range("J22").Value = 1052.5834
$ws.Range("K22").Value = 2295.2
$ws.Range("L22").Value = 1052.5834
$ws.Range("M22").Value = -2000.2
$ws.Range("N22").Value = -1642.5834

$ws.Range("H27").Value = 1418.0588
$ws.Range("I27").Value = 2295.2
$ws.Range("J27").Value = 1052.5834
$ws.Range("K27").Value = 2295.2
$ws.Range("L27").Value = 1052.5834
$ws.Range("M27").Value = -2188.2
$ws.Range("N27").Value = -1266.5834

$ws.Range("H93").Value = 29312.166
$ws.Range("I93").Value = 1673.909
$ws.Range("J93").Value = 333333
$ws.Range("K93").Value = 1673.909
$ws.Range("L93").Value = 333333
$ws.Range("M93").Value = -425.9090000000001
$ws.Range("N93").Value = -335829

$ws.Range("H95").Value = 26344
$ws.Range("J95").Value = 26344
$ws.Range("L95").Value = 26344
$ws.Range("N95").Value = -31836

$ws.Range("H122").Value = 15338.154
$ws.Range("I122").Value = 18239.6
$ws.Range("K122").Value = 54718.8
$ws.Range("M122").Value = -52268.8

$ws.Range("H126").Value = 68499.75
$ws.Range("I126").Value = 68499.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 205499.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -203029.25
$ws.Range("N126").ClearContents()

$ws.Range("H136").Value = 2703.3076
$ws.Range("I136").Value = 2703.3076
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8109.9228
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5559.9228
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3813.3333
$ws.Range("I81").Value = 3658.5
$ws.Range("K81").Value = 7317
$ws.Range("M81").Value = -6256

$ws.Range("H84").Value = 3813.3333
$ws.Range("I84").Value = 3658.5
$ws.Range("K84").Value = 36585
$ws.Range("M84").Value = -31281

$ws.Range("H107").Value = 31250726
$ws.Range("I107").Value = 645.8889
$ws.Range("K107").Value = 1937.6667
$ws.Range("M107").Value = -17.66670000000022

$ws.Range("H113").Value = 1007.7059
$ws.Range("I113").Value = 787.8
$ws.Range("J113").Value = 1321.8572
$ws.Range("K113").Value = 2363.4
$ws.Range("L113").Value = 3965.5716
$ws.Range("M113").Value = -193.3999999999996
$ws.Range("N113").Value = -8305.571599999999

$ws.Range("H126").Value = 5232
$ws.Range("I126").Value = 3964
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 11892
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -9422
$ws.Range("N126").Value = -24440

$ws.Range("H132").Value = 4954.8066
$ws.Range("I132").Value = 6199.5713
$ws.Range("J132").Value = 2340.8
$ws.Range("K132").Value = 18598.7139
$ws.Range("L132").Value = 7022.400000000001
$ws.Range("M132").Value = -16068.7139
$ws.Range("N132").Value = -12082.4
